$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Price" column (D): values like "59.953.14" or "1.00" must remain literal text,
# otherwise Excel auto-converts them to numbers (stripping thousands dots / trailing zeros).
# Force text format before assigning, then restore the default "Normal" style so no stray
# cell-style attribute is left behind.
$priceUpdates = [ordered]@{
    "D2" = "59.953.14"
    "D3" = "3.272.68"
    "D5" = "553.15"
    "D6" = "139.76"
    "D8" = "3.278.24"
    "D9" = "0.463"
    "D10" = "7.77"
    "D11" = "0.117"
    "D12" = "0.403"
    "D13" = "3.847.97"
    "D15" = "26.75"
    "D16" = "3.287.40"
    "D17" = "0.0000163"
    "D18" = "60.042.93"
    "D19" = "6.06"
    "D20" = "13.76"
    "D21" = "8.51"
    "D22" = "370.77"
    "D23" = "73.62"
    "D25" = "0.529"
    "D26" = "3.433.10"
    "D27" = "0.0000100"
    "D29" = "0.996"
    "D30" = "7.06"
    "D31" = "1.00"
    "D32" = "2.01"
    "D33" = "7.42"
    "D34" = "22.42"
    "D35" = "1.22"
    "D36" = "5.03"
    "D37" = "165.88"
    "D38" = "1.50"
    "D39" = "6.60"
    "D40" = "3.313.80"
    "D41" = "25.94"
    "D42" = "0.0722"
    "D43" = "41.56"
    "D44" = "0.743"
    "D45" = "4.07"
    "D47" = "1.55"
    "D49" = "2.328.83"
    "D50" = "6.32"
    "D51" = "21.07"
}
foreach ($ref in $priceUpdates.Keys) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $priceUpdates[$ref]
    $ws.Range($ref).Style = "Normal"
}

# "Volume(1h)" column (E): percentage strings already carry padding spaces, so Excel
# leaves them as text automatically.
$volumeUpdates = [ordered]@{
    "E2" = "  -3.43%  "
    "E3" = "  -4.35%  "
    "E4" = "  +0.02%  "
    "E5" = "  -4.37%  "
    "E6" = "  -8.30%  "
    "E7" = "  -0.07%  "
    "E8" = "  -4.18%  "
    "E9" = "  -4.09%  "
    "E10" = "  -3.51%  "
    "E11" = "  -5.53%  "
    "E12" = "  -3.32%  "
    "E13" = "  -3.99%  "
    "E14" = "  -0.30%  "
    "E15" = "  -6.85%  "
    "E16" = "  -3.82%  "
    "E17" = "  -5.12%  "
    "E18" = "  -3.39%  "
    "E19" = "  -6.76%  "
    "E20" = "  -5.52%  "
    "E21" = "  -4.94%  "
    "E22" = "  -3.15%  "
    "E23" = "  -2.09%  "
    "E24" = "  -0.04%  "
    "E25" = "  -7.23%  "
    "E26" = "  -3.65%  "
    "E27" = "  -10.56%  "
    "E28" = "  -6.08%  "
    "E29" = "  -0.21%  "
    "E30" = "  -8.10%  "
    "E31" = "  +0.07%  "
    "E32" = "  -5.13%  "
    "E33" = "  -6.15%  "
    "E34" = "  -3.45%  "
    "E35" = "  -8.91%  "
    "E36" = "  -7.99%  "
    "E37" = "  -1.62%  "
    "E38" = "  -7.06%  "
    "E39" = "  -4.76%  "
    "E40" = "  -4.06%  "
    "E41" = "  -16.58%  "
    "E42" = "  -7.88%  "
    "E43" = "  -2.80%  "
    "E44" = "  -4.56%  "
    "E45" = "  -7.52%  "
    "E46" = "  -6.21%  "
    "E47" = "  -7.42%  "
    "E48" = "  +0.10%  "
    "E49" = "  -8.39%  "
    "E50" = "  -7.99%  "
    "E51" = "  -6.68%  "
}
foreach ($ref in $volumeUpdates.Keys) {
    $ws.Range($ref).Value = $volumeUpdates[$ref]
}
